# ProjectEnquiry.xlsx touch-up:
#  - row 2 enquiry/reply text replaced with placeholder "eeee..." strings
#  - row 2 enquiry date/reply-date columns collapsed (F2 cleared, G2 given a new timestamp)
#  - row 2 ID / PROJECT ID re-entered as clean integers (was 1.0 / 3.0)
#  - a few column widths nudged
#  - selection left on G4 (as it was when the workbook was last saved)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 cell edits -------------------------------------------------
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = "eeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeee"
$ws.Range("E2").Value = "eeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeeee"
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = 45765.895975347223

# --- Column width tweaks ------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 12.43
$ws.Columns.Item(3).ColumnWidth = 8.96
$ws.Columns.Item(6).ColumnWidth = 11.96

# --- Selection ----------------------------------------------------------
$ws.Range("G4").Select()
